# Insert a new data row above the existing row 150, shifting rows
# 150:188 down to 151:189, and populate the new row 150 with the
# new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(150).Insert()

$ws.Cells.Item(150, 1).Value = 11
$ws.Cells.Item(150, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(150, 3).Value = "Bíobío"
$ws.Cells.Item(150, 4).Value = 44985
$ws.Cells.Item(150, 5).Value = 8
$ws.Cells.Item(150, 6).Value = 100112043
$ws.Cells.Item(150, 7).Value = "Pepino ensalada"
$ws.Cells.Item(150, 8).Value = "Sin especificar"
$ws.Cells.Item(150, 9).Value = "Primera"
$ws.Cells.Item(150, 10).Value = 220
$ws.Cells.Item(150, 11).Value = 6000
$ws.Cells.Item(150, 12).Value = 6500
$ws.Cells.Item(150, 13).Value = 6273
$ws.Cells.Item(150, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(150, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(150, 16).Value = 105
$ws.Cells.Item(150, 17).Value = 60
$ws.Cells.Item(150, 18).Value = "Hortaliza"
